$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows/cells that have no counterpart in the new layout ---
# (old A6:A11,A13 and B8:B11,B13 held data that no longer exists in the
#  rebuilt grid; clearing them drops both value and formatting so they
#  disappear from the saved sheetData entirely)
$ws.Range("A6:A11").Clear()
$ws.Range("A13").Clear()
$ws.Range("B8:B11").Clear()
$ws.Range("B13").Clear()

# --- New header row: A1:D1 ("Upper Men"/"Upper Women"/"Lower Men"/"Lower Women") ---
# A1/B1 already have the plain default font; C1/D1 are brand-new cells that
# also start on the default font, so one Bold flip covers all four.
$ws.Range("A1").Value = "Upper Men"
$ws.Range("B1").Value = "Upper Women"
$ws.Range("C1").Value = "Lower Men"
$ws.Range("D1").Value = "Lower Women"
$ws.Range("A1:D1").Font.Bold = $true

# --- Column A: Upper Men (keep A2/A3/A4's existing formatting, just retext) ---
$ws.Range("A2").Value = "Bert"
$ws.Range("A3").Value = "Ernie"
$ws.Range("A4").Value = "John Bobbitt"
$ws.Range("A5").Value = "Fuckin Frank"

# A2 keeps its old bold/Menlo font, only the color moves from green to
# automatic/theme text color.
$ws.Range("A2").Font.ThemeColor = 1

# --- Column B: Upper Women ---
$ws.Range("B2").Value = "Lorena Bobbitt"
$ws.Range("B3").Value = "Jenny"
$ws.Range("B4").Value = "Psycho Sarah"
$ws.Range("B5").Value = "Monica"
$ws.Range("B6").Value = "Amy"
$ws.Range("B7").Value = "Alexandria"

# --- Column C: Lower Men ---
$ws.Range("C2").Value = "Fuckin Frank"
$ws.Range("C3").Value = "Jeff"
$ws.Range("C4").Value = "Scott"
$ws.Range("C5").Value = "Aer"
$ws.Range("C6").Value = "James"
$ws.Range("C7").Value = "Michael"
$ws.Range("C8").Value = "Andrew"
$ws.Range("C9").Value = "Kevin"
$ws.Range("C10").Value = "Razi"

# --- Column D: Lower Women ---
$ws.Range("D2").Value = "Poppy"
$ws.Range("D3").Value = "Eleanor"
$ws.Range("D4").Value = "Julie"
$ws.Range("D5").Value = "Ellie"
$ws.Range("D6").Value = "Farrah"
$ws.Range("D7").Value = "Ava"

# --- Summary / check row ---
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 9
$ws.Range("D12").Value = 6
$ws.Range("E12").Formula = "=SUM(A12:D12)"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668
$ws.Columns.Item(2).ColumnWidth = 31.830729166666668
$ws.Columns.Item(3).ColumnWidth = 41.166666666666664
$ws.Columns.Item(4).ColumnWidth = 52.330729166666664

# --- Selection ---
$ws.Range("D13").Select()
